# Adds a new price-history snapshot column ("2026-02-03 21:23:48") right
# before the "nom" / "url_produit" columns, shifting those two columns one
# position to the right (EV->EW, EW->EX), matching the weekly LDLC price
# tracker's usual layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column EV (152) is where "nom" currently lives. Inserting a whole column
# there pushes "nom" -> EW and "url_produit" -> EX, and leaves a brand new,
# completely blank column EV in their place.
$ws.Range("EV1").EntireColumn.Insert()

# New snapshot timestamp header.
$ws.Range("EV1").Value = "2026-02-03 21:23:48"

# Carry forward the last known price (previously stored in column EU, the
# last snapshot before this one) into the new EV column, but only for rows
# whose tracking was still active at the previous snapshot (EU not blank).
# Rows where EU is already empty (product no longer tracked) stay blank.
$lastRow = $ws.Range("A1").End(-4121).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $prev = $ws.Cells.Item($r, 151).Value()
    if ($prev -ne $null -and $prev -ne "") {
        $ws.Cells.Item($r, 152).Value = $prev
    }
}
